$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")
$ws.Activate()

$map = @{
    "shirts/clsh-1"      = "shirts/clsh-1.jpeg"
    "shirts/clsh-2"      = "shirts/clsh-2.jpeg"
    "shirts/clsh-3"      = "shirts/clsh-3.jpeg"
    "outerwear/clou-1"   = "outerwear/clou-1.jpeg"
    "outerwear/clou-3"   = "outerwear/clou-3.jpeg"
    "outerwear/clou-2"   = "outerwear/clou-2.jpeg"
    "outerwear/clou-4"   = "outerwear/clou-4.jpeg"
    "outerwear/clou-5"   = "outerwear/clou-5.jpeg"
    "outerwear/clou-6"   = "outerwear/clou-6.jpeg"
    "outerwear/clou-7"   = "outerwear/clou-7.jpeg"
    "outerwear/clou-8"   = "outerwear/clou-8.jpeg"
    "outerwear/clou-9"   = "outerwear/clou-9.jpeg"
    "outerwear/clou-10"  = "outerwear/clou-10.jpeg"
    "outerwear/clou-11"  = "outerwear/clou-11.jpg"
    "outerwear/clou-12"  = "outerwear/clou-12.jpeg"
    "outerwear/clou-13"  = "outerwear/clou-13.jpg"
    "outerwear/clou-14"  = "outerwear/clou-14.jpg"
    "pants/clpa-1"       = "pants/clpa-1.jpg"
    "pants/clpa-2"       = "pants/clpa-2.jpg"
    "pants/clpa-3"       = "pants/clpa-3.jpg"
    "pants/clpa-4"       = "pants/clpa-4.jpg"
    "pants/clpa-5"       = "pants/clpa-5.jpg"
    "pants/clpa-6"       = "pants/clpa-6.jpg"
    "pants/clpa-7"       = "pants/clpa-7.jpg"
    "pants/clpa-8"       = "pants/clpa-8.jpg"
    "pants/clpa-9"       = "pants/clpa-9.jpg"
    "pants/clpa-10"      = "pants/clpa-10.jpg"
    "shirts/clsh-4"      = "shirts/clsh-4.jpg"
    "shirts/clsh-5"      = "shirts/clsh-5.jpg"
    "shirts/clsh-6"      = "shirts/clsh-6.jpg"
    "shirts/clsh-7"      = "shirts/clsh-7.jpg"
    "accessories/clac-6" = "accessories/clac-6.jpg"
    "accessories/clac-3" = "accessories/clac-3.jpg"
    "accessories/clac-5" = "accessories/clac-5.jpg"
    "accessories/clac-2" = "accessories/clac-2.jpg"
    "accessories/clac-4" = "accessories/clac-4.jpg"
    "accessories/clac-1" = "accessories/clac-1.jpg"
    "accessories/clac-7" = "accessories/clac-7.webp"
    "accessories/clac-8" = "accessories/clac-8.webp"
    "accessories/clac-9" = "accessories/clac-9.jpg"
    "gear/gear-1"        = "gear/gear-1.png"
    "gear/gear-2"        = "gear/gear-2.png"
    "gear/gear-3"        = "gear/gear-3.png"
    "gear/gear-4"        = "gear/gear-4.jpg"
    "gear/gear-5"        = "gear/gear-5.jpg"
    "gear/gear-6"        = "gear/gear-6.jpg"
    "gear/gear-7"        = "gear/gear-7.jpg"
}

# Apply the renames in the same order the product data was originally
# authored (grouped by category), so the workbook's shared-string table
# is rebuilt in the same sequence as the source edit.
$rowOrder = @(38,39,40,41,33,34,35,36,37,42,43,44,45,46,47,48,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,2,3,4,29,30,31,32)

foreach ($r in $rowOrder) {
    $cell = $ws.Cells.Item($r, 3)
    $old = $cell.Value2
    if ($map.ContainsKey($old)) {
        $cell.Value2 = $map[$old]
    }
}

$ws.Range("D10").Select()

